$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1697.3478
$ws.Range("I28").Value = 466.14285
$ws.Range("K28").Value = 466.14285
$ws.Range("M28").Value = 18.85714999999999

# Row 43
$ws.Range("H43").Value = 3928.1428
$ws.Range("I43").Value = 3599.4
$ws.Range("J43").Value = 4750
$ws.Range("K43").Value = 3599.4
$ws.Range("L43").Value = 4750
$ws.Range("M43").Value = -3530.4
$ws.Range("N43").Value = -4888

# Row 129
$ws.Range("H129").Value = 4066.3333
$ws.Range("J129").Value = 5499.5
$ws.Range("L129").Value = 16498.5
$ws.Range("N129").Value = -26498.5

# Row 137
$ws.Range("H137").Value = 2929.3333
$ws.Range("I137").Value = 4774.5
$ws.Range("K137").Value = 14323.5
$ws.Range("M137").Value = -11773.5

# Row 138
$ws.Range("H138").Value = 10105932
$ws.Range("J138").Value = 15879867
$ws.Range("L138").Value = 47639601
$ws.Range("N138").Value = -47649881

# Row 141
$ws.Range("H141").Value = 1496.3684
$ws.Range("I141").Value = 1496.3684
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4489.1052
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 690.8948
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 1000000000
$ws.Range("I8").Value = 1000000000
$ws.Range("K8").Value = 1000000000
$ws.Range("M8").Value = -999999856

# Row 12
$ws.Range("H12").Value = 901
$ws.Range("I12").Value = 1003
$ws.Range("J12").Value = 850
$ws.Range("K12").Value = 1003
$ws.Range("L12").Value = 850
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -1196

# Row 16
$ws.Range("H16").Value = 400164
$ws.Range("I16").Value = 1000000
$ws.Range("J16").Value = 273.33334
$ws.Range("K16").Value = 1000000
$ws.Range("L16").Value = 273.33334
$ws.Range("M16").Value = -999713
$ws.Range("N16").Value = -847.33334

# Row 36
$ws.Range("H36").Value = 4534.6
$ws.Range("I36").Value = 4534.6
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4534.6
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -4188.6
$ws.Range("N36").ClearContents()

# Row 133
$ws.Range("H133").Value = 70995
$ws.Range("I133").Value = 70995
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 70995
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -68465
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2557
$ws.Range("I134").Value = 2388.5518
$ws.Range("K134").Value = 7165.655400000001
$ws.Range("M134").Value = -4630.655400000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6482.8667
$ws.Range("J31").Value = 10263.333
$ws.Range("L31").Value = 10263.333
$ws.Range("N31").Value = -10853.333

# Row 34
$ws.Range("H34").Value = 6482.8667
$ws.Range("J34").Value = 10263.333
$ws.Range("L34").Value = 10263.333
$ws.Range("N34").Value = -10667.333

# Row 94
$ws.Range("H94").Value = 1463.45
$ws.Range("I94").Value = 1289
$ws.Range("J94").Value = 1557.3846
$ws.Range("K94").Value = 1289
$ws.Range("L94").Value = 1557.3846
$ws.Range("M94").Value = -838
$ws.Range("N94").Value = -2459.3846

# Row 99
$ws.Range("H99").Value = 8155.9473
$ws.Range("I99").Value = 8556.823
$ws.Range("K99").Value = 8556.823
$ws.Range("M99").Value = -7058.823

# Row 126
$ws.Range("H126").Value = 8155.9473
$ws.Range("I126").Value = 8556.823
$ws.Range("K126").Value = 25670.469
$ws.Range("M126").Value = -23200.469

# Row 140
$ws.Range("H140").Value = 200000
$ws.Range("J140").Value = 200000
$ws.Range("L140").Value = 200000
$ws.Range("N140").Value = -210360

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1112.7307
$ws.Range("I5").Value = 762.7143
$ws.Range("J5").Value = 1521.0834
$ws.Range("K5").Value = 2288.1429
$ws.Range("L5").Value = 4563.2502
$ws.Range("M5").Value = -2176.1429
$ws.Range("N5").Value = -4787.2502

# Row 7
$ws.Range("H7").Value = 417.36365
$ws.Range("J7").Value = 162.5
$ws.Range("L7").Value = 487.5
$ws.Range("N7").Value = -711.5

# Row 129
$ws.Range("H129").Value = 4038.5356
$ws.Range("J129").Value = 3655.7144
$ws.Range("L129").Value = 10967.1432
$ws.Range("N129").Value = -20967.1432

# Row 130
$ws.Range("H130").Value = 3143.3333
$ws.Range("I130").Value = 3030
$ws.Range("J130").Value = 3200
$ws.Range("K130").Value = 9090
$ws.Range("L130").Value = 9600
$ws.Range("M130").Value = -4070
$ws.Range("N130").Value = -19640

# Row 131
$ws.Range("H131").Value = 26004.912
$ws.Range("J131").Value = 4634.4517
$ws.Range("L131").Value = 13903.3551
$ws.Range("N131").Value = -23983.3551

# Row 134
$ws.Range("H134").Value = 6566.2856
$ws.Range("J134").Value = 13811.8
$ws.Range("L134").Value = 41435.39999999999
$ws.Range("N134").Value = -51575.39999999999

# Row 135
$ws.Range("H135").Value = 1112.7307
$ws.Range("I135").Value = 762.7143
$ws.Range("J135").Value = 1521.0834
$ws.Range("K135").Value = 6864.428699999999
$ws.Range("L135").Value = 13689.7506
$ws.Range("M135").Value = -4329.428699999999
$ws.Range("N135").Value = -18759.7506

# Row 136
$ws.Range("H136").Value = 2249.75
$ws.Range("I136").Value = 2249.75
$ws.Range("K136").Value = 6749.25
$ws.Range("M136").Value = -1649.25

# Row 139
$ws.Range("H139").Value = 3974.8333
$ws.Range("I139").Value = 3962.25
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 11886.75
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -6746.75
$ws.Range("N139").Value = -22280

# Row 141
$ws.Range("H141").Value = 7939.96
$ws.Range("I141").Value = 5852.4287
$ws.Range("J141").Value = 10596.818
$ws.Range("K141").Value = 17557.2861
$ws.Range("L141").Value = 31790.454
$ws.Range("M141").Value = -12377.2861
$ws.Range("N141").Value = -42150.454

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 55333
$ws.Range("I62").Value = 55333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 55333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -54647
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 55333
$ws.Range("I65").Value = 55333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 165999
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -162567
$ws.Range("N65").ClearContents()

# Row 80
$ws.Range("H80").Value = 4539.5835
$ws.Range("I80").Value = 4626.6665
$ws.Range("J80").Value = 4278.3335
$ws.Range("K80").Value = 4626.6665
$ws.Range("L80").Value = 4278.3335
$ws.Range("M80").Value = -3628.6665
$ws.Range("N80").Value = -6274.3335

# Row 83
$ws.Range("H83").Value = 4539.5835
$ws.Range("I83").Value = 4626.6665
$ws.Range("J83").Value = 4278.3335
$ws.Range("K83").Value = 23133.3325
$ws.Range("L83").Value = 21391.6675
$ws.Range("M83").Value = -18141.3325
$ws.Range("N83").Value = -31375.6675

$ws = $wb.Worksheets.Item("LTW")
# Row 50
$ws.Range("H50").Value = 35749.5
$ws.Range("J50").Value = 35999.668
$ws.Range("L50").Value = 35999.668
$ws.Range("N50").Value = -37273.668

# Row 55
$ws.Range("H55").Value = 715.8421
$ws.Range("J55").Value = 970
$ws.Range("L55").Value = 970
$ws.Range("N55").Value = -1316

# Row 136
$ws.Range("H136").Value = 3159.7646
$ws.Range("I136").Value = 1977.64
$ws.Range("K136").Value = 5932.92
$ws.Range("M136").Value = -3382.92

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 4517
$ws.Range("I132").Value = 4451.273
$ws.Range("K132").Value = 13353.819
$ws.Range("M132").Value = -10823.819

# Row 136
$ws.Range("H136").Value = 1394.2894
$ws.Range("I136").Value = 1126.8695
$ws.Range("K136").Value = 3380.6085
$ws.Range("M136").Value = -830.6085000000003
